$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row get reformatted as text and overwritten with data (bug:
# the "id"/"customerId"/"TotalAmount" headers are gone), then a few more
# rows of the same text-formatted id/customerId values are appended below.
$ws.Range("A1:C1").NumberFormat = "@"
$ws.Range("A2:B4").NumberFormat = "@"

$ws.Range("A1").Value = "1"
$ws.Range("B1").Value = "0"
$ws.Range("C1").Value = "1"

$ws.Range("A2").Value = "2"
$ws.Range("B2").Value = "0"

$ws.Range("A3").Value = "1"
$ws.Range("B3").Value = "0"

$ws.Range("A4").Value = "2"
$ws.Range("B4").Value = "0"
